# ran resolve and classify+summarise steps after changes to mapping file
$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet: clear column C data, zero out column B data ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2").Value = 0
$wsRange.Range("C2").ClearContents()
$wsRange.Range("C3").ClearContents()
$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()
$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()
$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()
$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# --- "Species qualification" sheet: Range Analysis count -> 0 ---
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("B5").Value = 0

# --- "High Priority break-up" sheet: add New High Species columns ---
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")
$wsBreakup.Range("D2").Value = 4
$wsBreakup.Range("E2").Value = 100
